$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new Natural Language / LTL Specification pair as row 51.
$ws.Range("A51").Value = "1. Proceed straight for 500 meters.`n2. Turn right at the traffic light.`n3. Continue straight for 800 meters.`n4. When you see the university campus, turn left.`n5. Proceed for another 1 kilometer to reach your destination."
$ws.Range("B51").Value = "G((goStraight -> (F[0,500] prepareRight)) & (makeRight -> (F[0,800] prepareLeft)) & (makeLeft -> (F[0,1000] reachDestination)))"

# Re-point all the existing data rows (2-50) at the already-identical
# "font1/no-fill/wrap" style so the redundant duplicate cellXf collapses away.
$r = $ws.Range("A2:B50")
$r.Font.Name = "Times New Roman"

# Leave the selection where the user would land after typing the new row.
$ws.Range("B52").Select()
